# Update "想去人数" (want-to-go count) figures for the two sheets that
# contain the full exhibition data: "展览" and "全部类型" (kept in sync).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 134
    $ws.Range("F3").Value = 1686
    $ws.Range("F6").Value = 461
    $ws.Range("F9").Value = 604
}
